$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.244664430618286
$ws.Range("B1").Value = 2.59766697883606
$ws.Range("C1").Value = 7.982121467590332
$ws.Range("D1").Value = 2.132286787033081
$ws.Range("E1").Value = 1.139056801795959
